$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -6.568799999999993
$ws.Range("D3").Value = -7.64279999999999
$ws.Range("D5").Value = -8.135399999999997
$ws.Range("E7").Value = 12.1158
$ws.Range("C9").Value = -11.77360000000001
$ws.Range("E9").Value = 14.12600000000002
$ws.Range("D11").Value = -8.391500000000002
$ws.Range("D12").Value = -8.365500000000004
$ws.Range("C13").Value = -12.4706
$ws.Range("C16").Value = -11.6786
$ws.Range("C18").Value = -14.11589999999999
$ws.Range("C20").Value = -13.75669999999999
$ws.Range("D21").Value = -7.758700000000004
$ws.Range("E21").Value = 13.22450000000001

$wb.Save()
